# Generate Report for Handoff
$wb = $excel.ActiveWorkbook

# Target column width (OOXML <col width="..."/>) is 17.2159881591797 chars
# (down from 29.9777047293527 - columns narrowed once the "Status" text
# shrank from "Handed back: in sync with en-US" to "Ready for handoff").
# This engine's ColumnWidth COM setter quantizes to whole pixels, so
# 16.3333333333333 is the closest input that lands nearest the target.
$newColWidth = 16.3333333333333

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 11:10:05"
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-09-02 11:09:56"
$wsZh.Columns.Item(3).ColumnWidth = $newColWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-09-02 11:10:05"
$wsDe.Columns.Item(3).ColumnWidth = $newColWidth
